$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 5.230988666666666
$ws.Range("H2").Value = 15.692966
$ws.Range("I2").Value = 0.2129406655351238
$ws.Range("J2").Value = 0.2129406655351238
$ws.Range("M2").Value = 1.101288666666667
$ws.Range("N2").Value = 3.303866
$ws.Range("O2").Value = 0.2351204158555016
$ws.Range("P2").Value = 0.2351204158555016
$ws.Range("Q2").Value = 5.760828534061778
$ws.Range("R2").Value = 51.847456806556
$ws.Range("S2").Value = 0.05006669783316557
$ws.Range("T2").Value = 0.05006669783316556
# Row 3
$ws.Range("G3").Value = 5.230988666666666
$ws.Range("H3").Value = 15.692966
$ws.Range("I3").Value = 0.2129406655351238
$ws.Range("J3").Value = 0.2129406655351238
$ws.Range("O3").Value = 0.3619536680130181
$ws.Range("P3").Value = 0.361953668013018
$ws.Range("Q3").Value = 8.868447306503553
$ws.Range("R3").Value = 79.81602575853198
$ws.Range("S3").Value = 0.0770746549595713
$ws.Range("T3").Value = 0.07707465495957129
# Row 4
$ws.Range("G4").Value = 5.230988666666666
$ws.Range("H4").Value = 15.692966
$ws.Range("I4").Value = 0.2129406655351238
$ws.Range("J4").Value = 0.2129406655351238
$ws.Range("M4").Value = 0.9788929999999999
$ws.Range("N4").Value = 2.936679
$ws.Range("O4").Value = 0.2089894649825745
$ws.Range("P4").Value = 0.2089894649825745
$ws.Range("Q4").Value = 5.120578188879332
$ws.Range("R4").Value = 46.08520369991399
$ws.Range("S4").Value = 0.04450235576321884
$ws.Range("T4").Value = 0.04450235576321884
# Row 5
$ws.Range("G5").Value = 5.230988666666666
$ws.Range("H5").Value = 15.692966
$ws.Range("I5").Value = 0.2129406655351238
$ws.Range("J5").Value = 0.2129406655351238
$ws.Range("M5").Value = 0.9083856666666668
$ws.Range("N5").Value = 2.725157
$ws.Range("O5").Value = 0.1939364511489059
$ws.Range("P5").Value = 0.1939364511489059
$ws.Range("Q5").Value = 4.751755127295778
$ws.Range("R5").Value = 42.765796145662
$ws.Range("S5").Value = 0.04129695697916803
$ws.Range("T5").Value = 0.04129695697916803
# Row 6
$ws.Range("I6").Value = 0.345577477529236
$ws.Range("J6").Value = 0.3455774775292359
$ws.Range("M6").Value = 1.101288666666667
$ws.Range("N6").Value = 3.303866
$ws.Range("O6").Value = 0.2351204158555016
$ws.Range("P6").Value = 0.2351204158555016
$ws.Range("Q6").Value = 9.349142345716666
$ws.Range("R6").Value = 84.14228111145
$ws.Range("S6").Value = 0.08125232022696921
$ws.Range("T6").Value = 0.0812523202269692
# Row 7
$ws.Range("I7").Value = 0.345577477529236
$ws.Range("J7").Value = 0.3455774775292359
$ws.Range("O7").Value = 0.3619536680130181
$ws.Range("P7").Value = 0.361953668013018
$ws.Range("S7").Value = 0.1250830355743933
$ws.Range("T7").Value = 0.1250830355743933
# Row 8
$ws.Range("I8").Value = 0.345577477529236
$ws.Range("J8").Value = 0.3455774775292359
$ws.Range("M8").Value = 0.9788929999999999
$ws.Range("N8").Value = 2.936679
$ws.Range("O8").Value = 0.2089894649825745
$ws.Range("P8").Value = 0.2089894649825745
$ws.Range("Q8").Value = 8.310091872574999
$ws.Range("R8").Value = 74.79082685317499
$ws.Range("S8").Value = 0.07222205213886267
$ws.Range("T8").Value = 0.07222205213886265
# Row 9
$ws.Range("I9").Value = 0.345577477529236
$ws.Range("J9").Value = 0.3455774775292359
$ws.Range("M9").Value = 0.9083856666666668
$ws.Range("N9").Value = 2.725157
$ws.Range("O9").Value = 0.1939364511489059
$ws.Range("P9").Value = 0.1939364511489059
$ws.Range("Q9").Value = 7.711535730391667
$ws.Range("R9").Value = 69.40382157352501
$ws.Range("S9").Value = 0.0670200695890108
$ws.Range("T9").Value = 0.06702006958901077
# Row 10
$ws.Range("G10").Value = 8.418577333333333
$ws.Range("H10").Value = 25.255732
$ws.Range("I10").Value = 0.3426995496362334
$ws.Range("J10").Value = 0.3426995496362334
$ws.Range("M10").Value = 1.101288666666667
$ws.Range("N10").Value = 3.303866
$ws.Range("O10").Value = 0.2351204158555016
$ws.Range("P10").Value = 0.2351204158555016
$ws.Range("Q10").Value = 9.271283806656889
$ws.Range("R10").Value = 83.44155425991201
$ws.Range("S10").Value = 0.0805756606239643
$ws.Range("T10").Value = 0.08057566062396429
# Row 11
$ws.Range("G11").Value = 8.418577333333333
$ws.Range("H11").Value = 25.255732
$ws.Range("I11").Value = 0.3426995496362334
$ws.Range("J11").Value = 0.3426995496362334
$ws.Range("O11").Value = 0.3619536680130181
$ws.Range("P11").Value = 0.361953668013018
$ws.Range("Q11").Value = 14.27258100407378
$ws.Range("R11").Value = 128.453229036664
$ws.Range("S11").Value = 0.124041359017244
$ws.Range("T11").Value = 0.124041359017244
# Row 12
$ws.Range("G12").Value = 8.418577333333333
$ws.Range("H12").Value = 25.255732
$ws.Range("I12").Value = 0.3426995496362334
$ws.Range("J12").Value = 0.3426995496362334
$ws.Range("M12").Value = 0.9788929999999999
$ws.Range("N12").Value = 2.936679
$ws.Range("O12").Value = 0.2089894649825745
$ws.Range("P12").Value = 0.2089894649825745
$ws.Range("Q12").Value = 8.240886421558665
$ws.Range("R12").Value = 74.167977794028
$ws.Range("S12").Value = 0.07162059552824564
$ws.Range("T12").Value = 0.07162059552824564
# Row 13
$ws.Range("G13").Value = 8.418577333333333
$ws.Range("H13").Value = 25.255732
$ws.Range("I13").Value = 0.3426995496362334
$ws.Range("J13").Value = 0.3426995496362334
$ws.Range("M13").Value = 0.9083856666666668
$ws.Range("N13").Value = 2.725157
$ws.Range("O13").Value = 0.1939364511489059
$ws.Range("P13").Value = 0.1939364511489059
$ws.Range("Q13").Value = 7.64731498332489
$ws.Range("R13").Value = 68.82583484992401
$ws.Range("S13").Value = 0.06646193446677943
$ws.Range("T13").Value = 0.06646193446677942
# Row 14
$ws.Range("G14").Value = 2.426634333333333
$ws.Range("H14").Value = 7.279902999999999
$ws.Range("I14").Value = 0.09878230729940689
$ws.Range("J14").Value = 0.09878230729940687
$ws.Range("M14").Value = 1.101288666666667
$ws.Range("N14").Value = 3.303866
$ws.Range("O14").Value = 0.2351204158555016
$ws.Range("P14").Value = 0.2351204158555016
$ws.Range("Q14").Value = 2.672424889444222
$ws.Range("R14").Value = 24.051824004998
$ws.Range("S14").Value = 0.0232257371714025
$ws.Range("T14").Value = 0.02322573717140249
# Row 15
$ws.Range("G15").Value = 2.426634333333333
$ws.Range("H15").Value = 7.279902999999999
$ws.Range("I15").Value = 0.09878230729940689
$ws.Range("J15").Value = 0.09878230729940687
$ws.Range("O15").Value = 0.3619536680130181
$ws.Range("P15").Value = 0.361953668013018
$ws.Range("Q15").Value = 4.114036578678443
$ws.Range("R15").Value = 37.02632920810599
$ws.Range("S15").Value = 0.03575461846180945
$ws.Range("T15").Value = 0.03575461846180945
# Row 16
$ws.Range("G16").Value = 2.426634333333333
$ws.Range("H16").Value = 7.279902999999999
$ws.Range("I16").Value = 0.09878230729940689
$ws.Range("J16").Value = 0.09878230729940687
$ws.Range("M16").Value = 0.9788929999999999
$ws.Range("N16").Value = 2.936679
$ws.Range("O16").Value = 0.2089894649825745
$ws.Range("P16").Value = 0.2089894649825745
$ws.Range("Q16").Value = 2.375415362459666
$ws.Range("R16").Value = 21.378738262137
$ws.Range("S16").Value = 0.0206444615522473
$ws.Range("T16").Value = 0.0206444615522473
# Row 17
$ws.Range("G17").Value = 2.426634333333333
$ws.Range("H17").Value = 7.279902999999999
$ws.Range("I17").Value = 0.09878230729940689
$ws.Range("J17").Value = 0.09878230729940687
$ws.Range("M17").Value = 0.9083856666666668
$ws.Range("N17").Value = 2.725157
$ws.Range("O17").Value = 0.1939364511489059
$ws.Range("P17").Value = 0.1939364511489059
$ws.Range("Q17").Value = 2.204319846641222
$ws.Range("R17").Value = 19.838878619771
$ws.Range("S17").Value = 0.01915749011394764
$ws.Range("T17").Value = 0.01915749011394763
